$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (unchanged text, but ensure it's set)
$ws.Range("A1").Value = "PERTEMUAN"
$ws.Range("B1").Value = "MATERI"
$ws.Range("C1").Value = "DOSEN"
$ws.Range("D1").Value = "MATAKULIAH"

# Set the MATAKULIAH column first so the new "Kalkulus I" shared string is
# created right after the header strings (matches author's edit order).
$ws.Range("D2").Value = "Kalkulus I"
$ws.Range("D3").Value = "Kalkulus I"
$ws.Range("D4").Value = "Kalkulus I"

# New MATERI values, per pertemuan, in row order 3, 4, 2 so the shared
# string table is built up in the same sequence as the source edit.
$ws.Range("B3").Value = "Rumus luas persegi"
$ws.Range("B4").Value = "Rumus keliling lingkaran"
$ws.Range("B2").Value = "Rumus luas segitiga"

# PERTEMUAN numbers renumbered 1, 2, 3
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# DOSEN stays "Andi A" for every row
$ws.Range("C2").Value = "Andi A"
$ws.Range("C3").Value = "Andi A"
$ws.Range("C4").Value = "Andi A"

# Move the active selection to B2
[void]$ws.Range("B2").Select()
